$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Restricciones_del_lider (MIU_value rows) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "2.09 - x"
Set-TextValue $ws2.Range("B2") "-3.09"
Set-TextValue $ws2.Range("D2") "0.86"

$ws2.Range("A3").Value = "-2.09 + x"
Set-TextValue $ws2.Range("B3") "1.0899999999999999"
Set-TextValue $ws2.Range("D3") "0.62"

$ws2.Range("A4").Value = "41.02289999999999 + x - y - 9(x^2)"
Set-TextValue $ws2.Range("B4") "-40.02289999999999"
Set-TextValue $ws2.Range("D4") "0.58"

# --- Restricciones_del_follower (Lambda/Beta/Gamma rows) ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "22.9596 - 12.084y + (-0.5 + x)*(y^2)"
Set-TextValue $ws3.Range("B2") "-22.9596"
Set-TextValue $ws3.Range("D2") "0.69"
Set-TextValue $ws3.Range("E2") "7.0"

$ws3.Range("A3").Value = "-1.8619999999999999 + 0.49y"
Set-TextValue $ws3.Range("B3") "0.8619999999999999"
Set-TextValue $ws3.Range("D3") "0.65"
Set-TextValue $ws3.Range("E3") "8.5"

Set-TextValue $ws3.Range("A4") "-2"
Set-TextValue $ws3.Range("B4") "-1"
Set-TextValue $ws3.Range("D4") "0.32"
Set-TextValue $ws3.Range("E4") "9.9"
Set-TextValue $ws3.Range("F4") "0"

# --- Punto_modificado (x, y point) ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "2.09"
Set-TextValue $ws4.Range("B2") "3.8"

# --- Vector_bf ---
# NOTE: sheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) lookup is case-insensitive, so both names would
# resolve to the same sheet. Use the 1-based positional index instead.
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-1.318499999999999"

# --- Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-83.27484000000001"
Set-TextValue $ws6.Range("A3") "-4.704041999999987"

# --- Vector_Alpha ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 0.51
